# Denied file upload for objects with DELETE and DELETED status.
#
# The workbook has three pairs of rules (Complaint, Case File, Task) that
# grant "uploadOrReplaceFile" access while an object is not CLOSED, and deny
# it once the object becomes CLOSED. This change extends both the "grant"
# and "deny" conditions so DELETE / DELETED statuses are treated the same
# way as CLOSED.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$grantCondition = "status != 'CLOSED' && status != 'DELETE' && status != 'DELETED'"
$denyCondition  = "status == 'CLOSED' || status == 'DELETE' || status == 'DELETED'"

# Row 26/27 -> Complaint "Only participants can upload or replace files" /
#              "Participants cannot upload or replace files"
# Row 41/42 -> Case File "Only participants can version files" /
#              "Participants cannot version files"
# Row 50/51 -> Task "Only participants can version files" /
#              "Participants cannot version files"
$grantRows = 26, 41, 50
$denyRows  = 27, 42, 51

foreach ($row in $grantRows) {
    $cell = $ws.Range("D$row")
    $cell.Value = $grantCondition
    $cell.WrapText = $true
}

foreach ($row in $denyRows) {
    $cell = $ws.Range("D$row")
    $cell.Value = $denyCondition
    $cell.WrapText = $true
}

# Reflect where the author ended up after making the edit.
[void]$ws.Range("D51").Select()
